$wb = $excel.ActiveWorkbook

$dateFmt = "yyyy-mm-dd HH:mm:ss"

# -----------------------------------------------------------------
# Overview sheet
# -----------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")
$lo = $ws.ListObjects.Item(1)

# Update existing row (file renamed + refreshed timestamp)
$ws.Range("A2").Value = "5d24a5f0-095a-4f6c-ba13-767e5b80d782.md"
$ws.Range("G2").Value = "2016-08-31 19:11:18"

# Add the new row for the second handed-back file
$lo.ListRows.Add() | Out-Null
$ws.Range("A3").Value = "84cb22ab-5658-4cb6-b7d7-fbe251bee46d.md"
$ws.Range("C3").Value = ".md"
$ws.Range("E3").Value = "Handed back: in sync with en-US"
$ws.Range("F3").Value = "Handed back: in sync with en-US"
$ws.Range("G3").Value = "2016-08-31 19:11:18"
$ws.Range("G3").NumberFormat = $dateFmt

# Rebuild the hyperlinks for column B so relationship ids stay sequential
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9109043c922aa5c25542bc756f1e19be63578c2b/e2e/5d24a5f0-095a-4f6c-ba13-767e5b80d782.md", "", "", "e2e\5d24a5f0-095a-4f6c-ba13-767e5b80d782.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9109043c922aa5c25542bc756f1e19be63578c2b/e2e/84cb22ab-5658-4cb6-b7d7-fbe251bee46d.md", "", "", "e2e\84cb22ab-5658-4cb6-b7d7-fbe251bee46d.md") | Out-Null

# -----------------------------------------------------------------
# zh-cn sheet
# -----------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")
$lo = $ws.ListObjects.Item(1)

# Update existing row (file renamed + refreshed hashes/timestamps)
$ws.Range("A2").Value = "5d24a5f0-095a-4f6c-ba13-767e5b80d782.md"
$ws.Range("G2").Value = "5d24a5f0-095a-4f6c-ba13-767e5b80d782.b0c7a6fa9b68ec5f8a10893d370c0d2f11b82a62.zh-cn.xlf"
$ws.Range("H2").Value = "2016-08-31 19:11:00"
$ws.Range("I2").Value = "5d24a5f0-095a-4f6c-ba13-767e5b80d782.md"
$ws.Range("J2").Value = "5d24a5f0-095a-4f6c-ba13-767e5b80d782.b0c7a6fa9b68ec5f8a10893d370c0d2f11b82a62.zh-cn.xlf"
$ws.Range("K2").Value = "2016-08-31 19:11:36"

# Add the new row for the second handed-back file
$lo.ListRows.Add() | Out-Null
$ws.Range("A3").Value = "84cb22ab-5658-4cb6-b7d7-fbe251bee46d.md"
$ws.Range("B3").Value = ".md"
$ws.Range("C3").Value = "Handed back: in sync with en-US"
$ws.Range("D3").Value = "e2e"
$ws.Range("E3").Value = "ht"
$ws.Range("F3").Value = "True"
$ws.Range("G3").Value = "84cb22ab-5658-4cb6-b7d7-fbe251bee46d.313c14a06a567d9a49dcb727d9d10f26a6f8a805.zh-cn.xlf"
$ws.Range("H3").Value = "2016-08-31 19:11:00"
$ws.Range("H3").NumberFormat = $dateFmt
$ws.Range("I3").Value = "84cb22ab-5658-4cb6-b7d7-fbe251bee46d.md"
$ws.Range("J3").Value = "84cb22ab-5658-4cb6-b7d7-fbe251bee46d.313c14a06a567d9a49dcb727d9d10f26a6f8a805.zh-cn.xlf"
$ws.Range("K3").Value = "2016-08-31 19:11:36"
$ws.Range("K3").NumberFormat = $dateFmt
$ws.Range("L3").Value = ""
$ws.Range("M3").Value = "True"
$ws.Range("N3").Value = ""
$ws.Range("O3").Value = "False"
$ws.Range("P3").Value = ""

# Rebuild the hyperlinks for columns A and I so relationship ids stay sequential
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9109043c922aa5c25542bc756f1e19be63578c2b/e2e/5d24a5f0-095a-4f6c-ba13-767e5b80d782.md", "", "", "5d24a5f0-095a-4f6c-ba13-767e5b80d782.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/a4e35a71d5e5abde4491f0e5b1bdb0679bcb9ebd/e2e/5d24a5f0-095a-4f6c-ba13-767e5b80d782.md", "", "", "5d24a5f0-095a-4f6c-ba13-767e5b80d782.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9109043c922aa5c25542bc756f1e19be63578c2b/e2e/84cb22ab-5658-4cb6-b7d7-fbe251bee46d.md", "", "", "84cb22ab-5658-4cb6-b7d7-fbe251bee46d.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/a4e35a71d5e5abde4491f0e5b1bdb0679bcb9ebd/e2e/84cb22ab-5658-4cb6-b7d7-fbe251bee46d.md", "", "", "84cb22ab-5658-4cb6-b7d7-fbe251bee46d.md") | Out-Null

# -----------------------------------------------------------------
# de-de sheet
# -----------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")
$lo = $ws.ListObjects.Item(1)

# Update existing row (file renamed + refreshed hashes/timestamps)
$ws.Range("A2").Value = "5d24a5f0-095a-4f6c-ba13-767e5b80d782.md"
$ws.Range("G2").Value = "5d24a5f0-095a-4f6c-ba13-767e5b80d782.b0c7a6fa9b68ec5f8a10893d370c0d2f11b82a62.de-de.xlf"
$ws.Range("H2").Value = "2016-08-31 19:11:18"
$ws.Range("I2").Value = "5d24a5f0-095a-4f6c-ba13-767e5b80d782.md"
$ws.Range("J2").Value = "5d24a5f0-095a-4f6c-ba13-767e5b80d782.b0c7a6fa9b68ec5f8a10893d370c0d2f11b82a62.de-de.xlf"
$ws.Range("K2").Value = "2016-08-31 19:11:44"

# Add the new row for the second handed-back file
$lo.ListRows.Add() | Out-Null
$ws.Range("A3").Value = "84cb22ab-5658-4cb6-b7d7-fbe251bee46d.md"
$ws.Range("B3").Value = ".md"
$ws.Range("C3").Value = "Handed back: in sync with en-US"
$ws.Range("D3").Value = "e2e"
$ws.Range("E3").Value = "ht"
$ws.Range("F3").Value = "True"
$ws.Range("G3").Value = "84cb22ab-5658-4cb6-b7d7-fbe251bee46d.313c14a06a567d9a49dcb727d9d10f26a6f8a805.de-de.xlf"
$ws.Range("H3").Value = "2016-08-31 19:11:18"
$ws.Range("H3").NumberFormat = $dateFmt
$ws.Range("I3").Value = "84cb22ab-5658-4cb6-b7d7-fbe251bee46d.md"
$ws.Range("J3").Value = "84cb22ab-5658-4cb6-b7d7-fbe251bee46d.313c14a06a567d9a49dcb727d9d10f26a6f8a805.de-de.xlf"
$ws.Range("K3").Value = "2016-08-31 19:11:44"
$ws.Range("K3").NumberFormat = $dateFmt
$ws.Range("L3").Value = ""
$ws.Range("M3").Value = "True"
$ws.Range("N3").Value = ""
$ws.Range("O3").Value = "False"
$ws.Range("P3").Value = ""

# Rebuild the hyperlinks for columns A and I so relationship ids stay sequential
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9109043c922aa5c25542bc756f1e19be63578c2b/e2e/5d24a5f0-095a-4f6c-ba13-767e5b80d782.md", "", "", "5d24a5f0-095a-4f6c-ba13-767e5b80d782.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/2427b5bc894896146e715bbb88d42ee3a149ff70/e2e/5d24a5f0-095a-4f6c-ba13-767e5b80d782.md", "", "", "5d24a5f0-095a-4f6c-ba13-767e5b80d782.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9109043c922aa5c25542bc756f1e19be63578c2b/e2e/84cb22ab-5658-4cb6-b7d7-fbe251bee46d.md", "", "", "84cb22ab-5658-4cb6-b7d7-fbe251bee46d.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/2427b5bc894896146e715bbb88d42ee3a149ff70/e2e/84cb22ab-5658-4cb6-b7d7-fbe251bee46d.md", "", "", "84cb22ab-5658-4cb6-b7d7-fbe251bee46d.md") | Out-Null
